# [FEATURE] Add tablet homescreen and own profile
#
# Adds two new time-tracking entries to the "Arbeitsmatrix" sheet
# (rows 148 and 149), pushing the existing blank-separator row and the
# summary block (Summe / Budget / Kredits) further down. A second blank
# separator row is also introduced so the summary block keeps its usual
# spacing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 3 new rows right before the old blank separator / summary block.
# Excel copies row 147's (the row above the insertion point) formatting
# into all 3 new rows, which matches the styling the two new data rows
# need (s=21 for A/B/C/E, s=13 for D, s=2 for F/G, s=15 for I, s=19 for J/K).
$ws.Rows.Item(148).Resize(3).Insert()

# Row 148: new [TASK] entry "Mobile refinement"
$ws.Range("A148").Value = 22
$ws.Range("B148").Value = "Interface Design"
$ws.Range("C148").Value = "MockUps"
$ws.Range("D148").Value = "[TASK]"
$ws.Range("E148").Value = "Mobile refinement"
$ws.Range("F148").Value = 44491
$ws.Range("G148").Value = 44481
$ws.Range("I148").Formula = "=ROUNDUP(((SUM(K148-J148)*24*60/60)/0.25),0)*0.25"
$ws.Range("J148").Value = 0.45833333333333331
$ws.Range("K148").Value = 0.5

# Row 149: new [FEATURE] entry "Tablet Startseite und Profiluebersicht"
$ws.Range("A149").Value = 22
$ws.Range("B149").Value = "Interface Design"
$ws.Range("C149").Value = "MockUps"
$ws.Range("D149").Value = "[FEATURE]"
$ws.Range("E149").Value = "Tablet Startseite und Profilübersicht"
$ws.Range("F149").Value = 44491
$ws.Range("G149").Value = 44481
$ws.Range("I149").Formula = "=ROUNDUP(((SUM(K149-J149)*24*60/60)/0.25),0)*0.25"
$ws.Range("J149").Value = 0.5
$ws.Range("K149").Value = 0.66666666666666663

# Row 150 is the freshly-inserted row that must stay an empty separator
# row (like row 146/151): only D/F/G carry formatting, no other cells.
# Insert() stamped every column with row 147's style, so wipe it and
# restore just the D/F/G formatting from the existing separator row 146.
$ws.Range("A150:K150").Clear()
$ws.Range("D146").Copy()
$ws.Range("D150").PasteSpecial(-4122)
$ws.Range("F146:G146").Copy()
$ws.Range("F150:G150").PasteSpecial(-4122)

# Keep the selection / active cell in sync with the now-larger sheet.
$ws.Range("I155").Select() | Out-Null
